$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data and row reorderings
# Force text format on target cells first so numeric-looking strings
# (e.g. "31.50", "0.999") are preserved exactly as text, matching the source data.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "34.426.78"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +12.67%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.823.88"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +7.97%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "231.20"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +4.67%  "
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +4.42%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "31.50"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +2.22%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "45.76"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +3.38%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.283"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +6.29%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +8.18%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +3.22%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.081.39"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +7.64%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.798.17"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +6.93%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.647"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +4.33%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "34.406.70"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +12.44%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "10.29"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -4.39%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.35"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +7.75%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "70.19"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +6.17%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "260.73"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +5.24%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0752"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +4.45%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.57"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +3.77%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +2.07%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +1.27%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "161.05"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +2.31%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.82"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +5.43%  "
$ws.Range("B28").NumberFormat = "@"
$ws.Range("B28").Value = "Cosmos"
$ws.Range("C28").NumberFormat = "@"
$ws.Range("C28").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.17"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +5.67%  "
$ws.Range("B29").NumberFormat = "@"
$ws.Range("B29").Value = "Stellar"
$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.117"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +4.47%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.10%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.84"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +9.60%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0517"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +3.13%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +6.87%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.59"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +8.36%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.588.71"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +5.70%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +5.06%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +3.20%  "
$ws.Range("B38").NumberFormat = "@"
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").NumberFormat = "@"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0190"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +5.43%  "
$ws.Range("B39").NumberFormat = "@"
$ws.Range("B39").Value = "ImmutableX"
$ws.Range("C39").NumberFormat = "@"
$ws.Range("C39").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.633"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +7.75%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "85.11"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +7.19%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.89"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +6.72%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +1.19%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.919"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +7.51%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.16"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +6.83%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0521"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +3.70%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +4.38%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.974.01"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +7.88%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.75"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +5.62%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "53.46"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +2.10%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.999"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.12%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0₆0123"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +7.86%  "
